$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("P236:S236").Select()
$excel.ActiveWindow.ScrollRow = 217
